$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: rename label from YC_CTRS to YC_Default
$ws.Range("B2").Value2 = "YC_Default"

# Add new row 3: duplicate of row 2 but labeled YC_Shifted
$ws.Range("A2:M2").Copy($ws.Range("A3:M3"))

$ws.Range("A3").Value2 = "YieldCurve"
$ws.Range("B3").Value2 = "YC_Shifted"
$ws.Range("C3").Value2 = 42736
$ws.Range("D3").Value2 = "1Y"
$ws.Range("E3").Value2 = 0.0050000000000000001
$ws.Range("F3").Value2 = "2Y"
$ws.Range("G3").Value2 = 0.01
$ws.Range("H3").Value2 = "5Y"
$ws.Range("I3").Value2 = 0.02
$ws.Range("J3").Value2 = "10Y"
$ws.Range("K3").Value2 = 0.03
$ws.Range("L3").Value2 = "15Y"
$ws.Range("M3").Value2 = 0.04

# Update the selection to reflect the saved view state
$ws.Range("B9").Select()

Write-Output "done"
